$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting for numeric-looking price values so Excel
# keeps them as text (matching the original inline-string cells)
# instead of auto-converting to floating point numbers.
$textCells = @('D5', 'D7', 'D8', 'D9', 'D10', 'D11', 'D12', 'D14', 'D15', 'D16', 'D19', 'D20', 'D21', 'D22', 'D24', 'D25', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D49', 'D50', 'D51')
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '30.470.58'
$ws.Range('E2').Value = '  +0.21%  '
$ws.Range('D3').Value = '2.108.92'
$ws.Range('E3').Value = '  +1.77%  '
$ws.Range('E4').Value = '  -0.32%  '
$ws.Range('D5').Value = '336.51'
$ws.Range('E5').Value = '  +2.34%  '
$ws.Range('E6').Value = '  -0.49%  '
$ws.Range('D7').Value = '0.5239'
$ws.Range('E7').Value = '  +0.77%  '
$ws.Range('D8').Value = '0.4585'
$ws.Range('E8').Value = '  +5.58%  '
$ws.Range('D9').Value = '53.25'
$ws.Range('E9').Value = '  +15.80%  '
$ws.Range('D10').Value = '0.08927'
$ws.Range('E10').Value = '  +3.25%  '
$ws.Range('D11').Value = '1.180'
$ws.Range('E11').Value = '  +2.48%  '
$ws.Range('D12').Value = '24.52'
$ws.Range('E12').Value = '  +1.50%  '
$ws.Range('D13').Value = '2.094.07'
$ws.Range('E13').Value = '  +0.55%  '
$ws.Range('D14').Value = '6.806'
$ws.Range('E14').Value = '  +2.86%  '
$ws.Range('D15').Value = '8.000'
$ws.Range('E15').Value = '  +4.15%  '
$ws.Range('D16').Value = '96.58'
$ws.Range('E16').Value = '  +1.46%  '
$ws.Range('E17').Value = '  -0.32%  '
$ws.Range('E18').Value = '  +1.95%  '
$ws.Range('D19').Value = '0.06630'
$ws.Range('E19').Value = '  +0.14%  '
$ws.Range('D20').Value = '19.33'
$ws.Range('E20').Value = '  +3.39%  '
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').Value = '  -0.51%  '
$ws.Range('D22').Value = '6.309'
$ws.Range('E22').Value = '  +1.28%  '
$ws.Range('D23').Value = '30.519.06'
$ws.Range('E23').Value = '  +0.24%  '
$ws.Range('D24').Value = '12.38'
$ws.Range('E24').Value = '  +1.22%  '
$ws.Range('D25').Value = '2.364'
$ws.Range('E25').Value = '  +2.23%  '
$ws.Range('D26').Value = '2.334.01'
$ws.Range('E26').Value = '  +0.21%  '
$ws.Range('D27').Value = '22.38'
$ws.Range('E27').Value = '  +1.14%  '
$ws.Range('D28').Value = '2.573'
$ws.Range('E28').Value = '  +2.53%  '
$ws.Range('D29').Value = '163.49'
$ws.Range('E29').Value = '  +1.17%  '
$ws.Range('D30').Value = '132.86'
$ws.Range('E30').Value = '  +1.74%  '
$ws.Range('D31').Value = '1.224'
$ws.Range('E31').Value = '  +4.38%  '
$ws.Range('D32').Value = '1.717'
$ws.Range('E32').Value = '  +14.87%  '
$ws.Range('E33').Value = '  +0.60%  '
$ws.Range('D34').Value = '6.204'
$ws.Range('E34').Value = '  +2.88%  '
$ws.Range('D35').Value = '3.927'
$ws.Range('E35').Value = '  +2.82%  '
$ws.Range('D36').Value = '10.52'
$ws.Range('E36').Value = '  +9.59%  '
$ws.Range('D37').Value = '0.02580'
$ws.Range('E37').Value = '  +1.02%  '
$ws.Range('D38').Value = '0.06843'
$ws.Range('E38').Value = '  +3.88%  '
$ws.Range('D39').Value = '5.560'
$ws.Range('E39').Value = '  +2.26%  '
$ws.Range('D40').Value = '12.88'
$ws.Range('E40').Value = '  +3.72%  '
$ws.Range('D41').Value = '0.2306'
$ws.Range('E41').Value = '  +3.17%  '
$ws.Range('D42').Value = '0.6910'
$ws.Range('E42').Value = '  +3.11%  '
$ws.Range('D43').Value = '1.248'
$ws.Range('E43').Value = '  +1.37%  '
$ws.Range('D44').Value = '2.356'
$ws.Range('E44').Value = '  +7.68%  '
$ws.Range('D45').Value = '1.001'
$ws.Range('E45').Value = '  -0.42%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = '0.6389'
$ws.Range('E46').Value = '  +1.54%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '13.97'
$ws.Range('E47').Value = '  +1.12%  '
$ws.Range('E48').Value = '  +1.29%  '
$ws.Range('D49').Value = '0.00000000353'
$ws.Range('E49').Value = '  +26.26%  '
$ws.Range('D50').Value = '1.249'
$ws.Range('E50').Value = '  +1.64%  '
$ws.Range('D51').Value = '83.79'
$ws.Range('E51').Value = '  +2.91%  '
